$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update List Price values (column C) that changed
$ws.Range("C5").Value = 293.78
$ws.Range("C6").Value = 299.49
$ws.Range("C8").Value = 51.24
$ws.Range("C9").Value = 52.56

# Update the Notes (column J) text that changed
$ws.Range("J8").Value = "2 listings found (excluding us). Price set to .3*First + .7*Second"

# Update the selection to A12
$ws.Range("A12").Select()
